$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the "Step, Dir to
# ESS." paragraph. After the edits below it ends up at the start of the
# "Aux relay outputs (24V)" paragraph (right after its new "+ " prefix),
# so drop it now and re-add it in the right spot once that text exists.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Add "+ " prefix to several list items that didn't have one yet.
$d.Content.Find.Execute("Brake relay", $true, $false, $false, $false, $false, $true, 1, $false, "+ Brake relay", 2)
$d.Content.Find.Execute("Limit switches circuit", $true, $false, $false, $false, $false, $true, 1, $false, "+ Limit switches circuit", 2)
$d.Content.Find.Execute("Connect to ESS", $true, $false, $false, $false, $false, $true, 1, $false, "+ Connect to ESS", 2)
$d.Content.Find.Execute("Step, Dir to ESS.", $true, $false, $false, $false, $false, $true, 1, $false, "+ Step, Dir, Enable to ESS.", 2)
$d.Content.Find.Execute("General purpose I/O", $true, $false, $false, $false, $false, $true, 1, $false, "+ General purpose I/O", 2)
$d.Content.Find.Execute("Aux relay outputs (24V)", $true, $false, $false, $false, $false, $true, 1, $false, "+ Aux relay outputs (24V)", 2)
$d.Content.Find.Execute("Earth to GND connection", $true, $false, $false, $false, $false, $true, 1, $false, "+ Earth to GND connection", 2)

# Re-add "_GoBack" right before "Aux relay outputs (24V)" (i.e. right
# after the "+ " that was just added in front of it).
$auxPara = $d.Paragraphs(10)
$bmPos = $auxPara.Range.Start + 2
$insertPoint = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $insertPoint)

# Fill in the final empty bullet item (it already exists as an empty
# ListParagraph at the end of the document -- just set its text instead
# of inserting a new paragraph).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "Add footprints wherever missing"
